$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I so the existing "Transfer price(per kg)"
# column (I) shifts to J, making room for the new "Distribution channel code"
# column.
$ws.Columns("I").Insert()

# New header in I1, matching the style of the other header cells (bold).
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I1").Font.Bold = $true

# New data values.
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Column width: new column I gets a manually-set (non best-fit) width; the
# shifted "Transfer price" column (now J) automatically keeps its original
# best-fit width from the Insert(), so it does not need to be touched.
$ws.Columns("I").ColumnWidth = 21.666666666666664

# Update the view selection state to match the post-edit selection.
$ws.Columns("I").Select()
